$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 3141
$ws.Range("I3").Value = 3226
$ws.Range("I4").Value = 764
$ws.Range("I5").Value = 299
$ws.Range("I6").Value = 3680
$ws.Range("I7").Value = 11110

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I3").Value = 26
$ws.Range("I4").Value = 15
$ws.Range("I7").Value = 121

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("I2").Value = 35
$ws.Range("I7").Value = 127

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 115
$ws.Range("I3").Value = 109
$ws.Range("I7").Value = 360

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I3").Value = 70
$ws.Range("I7").Value = 203

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I2").Value = 106
$ws.Range("I3").Value = 156
$ws.Range("I5").Value = 10
$ws.Range("I6").Value = 153
$ws.Range("I7").Value = 447

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I3").Value = 30
$ws.Range("I7").Value = 99

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("I2").Value = 36
$ws.Range("I7").Value = 97

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I3").Value = 63
$ws.Range("I6").Value = 81
$ws.Range("I7").Value = 246

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I4").Value = 44
$ws.Range("I5").Value = 36
$ws.Range("I7").Value = 375
$ws.Range("I8").Value = 702
$ws.Range("I19").Value = 297
$ws.Range("I20").Value = 284
$ws.Range("I21").Value = 61
$ws.Range("I23").Value = 100
$ws.Range("I27").Value = 96
$ws.Range("I29").Value = 728
$ws.Range("I31").Value = 99
$ws.Range("I33").Value = 502
$ws.Range("I36").Value = 151
$ws.Range("I37").Value = 360
$ws.Range("I42").Value = 391
$ws.Range("I45").Value = 21
$ws.Range("I47").Value = 80
$ws.Range("I48").Value = 129
$ws.Range("I49").Value = 83
$ws.Range("I50").Value = 51
$ws.Range("I52").Value = 238
$ws.Range("I54").Value = 247
$ws.Range("I55").Value = 121
$ws.Range("I59").Value = 21
$ws.Range("I63").Value = 44
$ws.Range("I64").Value = 102
$ws.Range("I65").Value = 246
$ws.Range("I67").Value = 447
$ws.Range("I69").Value = 27
$ws.Range("I72").Value = 38
$ws.Range("I73").Value = 93
$ws.Range("I78").Value = 152
$ws.Range("I79").Value = 283
$ws.Range("I80").Value = 37
$ws.Range("I83").Value = 226
$ws.Range("I84").Value = 97
$ws.Range("I85").Value = 512
$ws.Range("I86").Value = 65
$ws.Range("I89").Value = 121
$ws.Range("I90").Value = 135
$ws.Range("I93").Value = 62
$ws.Range("I95").Value = 176
$ws.Range("I96").Value = 127
$ws.Range("I97").Value = 86
$ws.Range("I99").Value = 203
$ws.Range("I101").Value = 11110

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("I2").Value = 79
$ws.Range("I6").Value = 43
$ws.Range("I7").Value = 226

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 61
$ws.Range("I7").Value = 176

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I6").Value = 163
$ws.Range("I7").Value = 502

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("I6").Value = 51
$ws.Range("I7").Value = 83

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("I6").Value = 126
$ws.Range("I7").Value = 247

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I2").Value = 222
$ws.Range("I3").Value = 255
$ws.Range("I4").Value = 29
$ws.Range("I6").Value = 195
$ws.Range("I7").Value = 728

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I2").Value = 119
$ws.Range("I7").Value = 297

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("I2").Value = 18
$ws.Range("I3").Value = 25
$ws.Range("I4").Value = 12
$ws.Range("I6").Value = 73
$ws.Range("I7").Value = 129

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 131
$ws.Range("I3").Value = 204
$ws.Range("I6").Value = 131
$ws.Range("I7").Value = 512

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I2").Value = 104
$ws.Range("I3").Value = 134
$ws.Range("I5").Value = 15
$ws.Range("I6").Value = 105
$ws.Range("I7").Value = 391

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I3").Value = 39
$ws.Range("I7").Value = 152

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I2").Value = 40
$ws.Range("I7").Value = 121

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("I2").Value = 28
$ws.Range("I6").Value = 29
$ws.Range("I7").Value = 100

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("I2").Value = 11
$ws.Range("I7").Value = 27

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("I6").Value = 49
$ws.Range("I7").Value = 61

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I3").Value = 89
$ws.Range("I5").Value = 10
$ws.Range("I7").Value = 283

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I3").Value = 33
$ws.Range("I7").Value = 102

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("I5").Value = 10
$ws.Range("I7").Value = 284

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("I3").Value = 47
$ws.Range("I7").Value = 151

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("I6").Value = 25
$ws.Range("I7").Value = 62

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I2").Value = 64
$ws.Range("I7").Value = 238

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I2").Value = 14
$ws.Range("I7").Value = 80

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I3").Value = 14
$ws.Range("I7").Value = 51

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 93

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("I2").Value = 9
$ws.Range("I7").Value = 21

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("I2").Value = 17
$ws.Range("I6").Value = 49
$ws.Range("I7").Value = 86

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 219
$ws.Range("I3").Value = 192
$ws.Range("I5").Value = 22
$ws.Range("I6").Value = 226
$ws.Range("I7").Value = 702

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("I2").Value = 6
$ws.Range("I7").Value = 36

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I3").Value = 18
$ws.Range("I7").Value = 96

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("I6").Value = 13
$ws.Range("I7").Value = 65

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I3").Value = 27
$ws.Range("I7").Value = 135

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("I6").Value = 20
$ws.Range("I7").Value = 38

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("I4").Value = 1
$ws.Range("I6").Value = 5
$ws.Range("I7").Value = 21

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("I2").Value = 7
$ws.Range("I6").Value = 18
$ws.Range("I7").Value = 37

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 129
$ws.Range("I3").Value = 110
$ws.Range("I7").Value = 375

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("I3").Value = 12
$ws.Range("I7").Value = 44
